$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update raw data values
$ws.Range("B2").Value = 30
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 12
$ws.Range("C7").Value = 339
$ws.Range("D7").Value = 87
$ws.Range("C13").Value = 1200
$ws.Range("D13").Value = 900

# Force recalculation so the SUBTOTAL/SUM formulas in row 14 pick up new totals
$excel.Calculate()

# Update the selected cell/view to B3
$ws.Range("B3").Select()
